$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 43
$ws.Range("I6").Value = 44.666668
$ws.Range("J6").Value = 38
$ws.Range("K6").Value = 134.000004
$ws.Range("L6").Value = 114
$ws.Range("M6").Value = -22.00000399999999
$ws.Range("N6").Value = -338
$ws.Range("H8").Value = 70.545456
$ws.Range("I8").Value = 70.545456
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 211.636368
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -72.636368
$ws.Range("N8").ClearContents()
$ws.Range("I9").Value = 111111300
$ws.Range("J9").Value = 1675.6666
$ws.Range("K9").Value = 111111300
$ws.Range("L9").Value = 1675.6666
$ws.Range("M9").Value = -111111131
$ws.Range("N9").Value = -2013.6666
$ws.Range("H33").Value = 335.6154
$ws.Range("I33").Value = 335.6154
$ws.Range("K33").Value = 335.6154
$ws.Range("M33").Value = -106.6154
$ws.Range("H52").Value = 9
$ws.Range("I52").Value = 9
$ws.Range("K52").Value = 27
$ws.Range("M52").Value = 133
$ws.Range("H113").Value = 3809.5789
$ws.Range("I113").Value = 5675
$ws.Range("J113").Value = 3312.1333
$ws.Range("K113").Value = 5675
$ws.Range("L113").Value = 3312.1333
$ws.Range("M113").Value = -2421
$ws.Range("N113").Value = -9820.1333
$ws.Range("H134").Value = 76616.91
$ws.Range("J134").Value = 76616.91
$ws.Range("L134").Value = 76616.91
$ws.Range("N134").Value = -86756.91
$ws.Range("H136").Value = 78935.42999999999
$ws.Range("J136").Value = 78935.42999999999
$ws.Range("L136").Value = 78935.42999999999
$ws.Range("N136").Value = -89135.42999999999
$ws.Range("H138").Value = 2970.9375
$ws.Range("I138").Value = 2839.1428
$ws.Range("J138").Value = 3073.4443
$ws.Range("K138").Value = 8517.428400000001
$ws.Range("L138").Value = 9220.332900000001
$ws.Range("M138").Value = -3377.428400000001
$ws.Range("N138").Value = -19500.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 27163
$ws.Range("J7").Value = 27163
$ws.Range("L7").Value = 27163
$ws.Range("N7").Value = -27391
$ws.Range("H32").Value = 2726.3562
$ws.Range("I32").Value = 2039.5605
$ws.Range("J32").Value = 9201.857
$ws.Range("K32").Value = 2039.5605
$ws.Range("L32").Value = 9201.857
$ws.Range("M32").Value = -1752.5605
$ws.Range("N32").Value = -9775.857
$ws.Range("H52").Value = 53620.285
$ws.Range("J52").Value = 53620.285
$ws.Range("L52").Value = 53620.285
$ws.Range("N52").Value = -54256.285
$ws.Range("H107").Value = 33803.855
$ws.Range("J107").Value = 33803.855
$ws.Range("L107").Value = 33803.855
$ws.Range("N107").Value = -41483.855
$ws.Range("H108").Value = 58233
$ws.Range("J108").Value = 58233
$ws.Range("L108").Value = 58233
$ws.Range("N108").Value = -65913
$ws.Range("H110").Value = 1138.2
$ws.Range("J110").Value = 2108.6667
$ws.Range("L110").Value = 2108.6667
$ws.Range("N110").Value = -6198.6667
$ws.Range("H118").Value = 49220
$ws.Range("J118").Value = 49220
$ws.Range("L118").Value = 49220
$ws.Range("N118").Value = -52534
$ws.Range("H135").Value = 102174
$ws.Range("J135").Value = 102174
$ws.Range("L135").Value = 102174
$ws.Range("N135").Value = -112314

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H19").Value = 1650
$ws.Range("I19").Value = 1650
$ws.Range("K19").Value = 1650
$ws.Range("M19").Value = -1477
$ws.Range("H52").Value = 99990
$ws.Range("J52").Value = 99990
$ws.Range("L52").Value = 99990
$ws.Range("N52").Value = -100516
$ws.Range("H55").Value = 57749
$ws.Range("J55").Value = 57749
$ws.Range("L55").Value = 57749
$ws.Range("N55").Value = -58295
$ws.Range("H99").Value = 1113767.9
$ws.Range("I99").Value = 36163.344
$ws.Range("J99").Value = 4586049
$ws.Range("K99").Value = 36163.344
$ws.Range("L99").Value = 4586049
$ws.Range("M99").Value = -34665.344
$ws.Range("N99").Value = -4589045
$ws.Range("H107").Value = 2883.6428
$ws.Range("I107").Value = 2670.6
$ws.Range("K107").Value = 2670.6
$ws.Range("M107").Value = -750.5999999999999
$ws.Range("H115").Value = 85278.71000000001
$ws.Range("J115").Value = 94990
$ws.Range("L115").Value = 94990
$ws.Range("N115").Value = -98124
$ws.Range("H121").Value = 99990
$ws.Range("J121").Value = 99990
$ws.Range("L121").Value = 99990
$ws.Range("N121").Value = -103484
$ws.Range("H127").Value = 57734
$ws.Range("J127").Value = 57734
$ws.Range("L127").Value = 57734
$ws.Range("N127").Value = -67654
$ws.Range("H132").Value = 71659.44500000001
$ws.Range("J132").Value = 77497.5
$ws.Range("L132").Value = 77497.5
$ws.Range("N132").Value = -87617.5
$ws.Range("H135").Value = 99448.57000000001
$ws.Range("J135").Value = 99448.57000000001
$ws.Range("L135").Value = 99448.57000000001
$ws.Range("N135").Value = -109588.57
$ws.Range("H138").Value = 79501.42999999999
$ws.Range("J138").Value = 79501.42999999999
$ws.Range("L138").Value = 79501.42999999999
$ws.Range("N138").Value = -89781.42999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1168.8064
$ws.Range("I58").Value = 949.46155
$ws.Range("J58").Value = 2309.4
$ws.Range("K58").Value = 949.46155
$ws.Range("L58").Value = 2309.4
$ws.Range("M58").Value = -746.46155
$ws.Range("N58").Value = -2715.4
$ws.Range("H118").Value = 62007.668
$ws.Range("J118").Value = 62007.668
$ws.Range("L118").Value = 62007.668
$ws.Range("N118").Value = -65321.668
$ws.Range("H119").Value = 63886.5
$ws.Range("J119").Value = 63886.5
$ws.Range("L119").Value = 63886.5
$ws.Range("N119").Value = -73562.5
$ws.Range("H136").Value = 1168.8064
$ws.Range("I136").Value = 949.46155
$ws.Range("J136").Value = 2309.4
$ws.Range("K136").Value = 2848.38465
$ws.Range("L136").Value = 6928.200000000001
$ws.Range("M136").Value = -298.38465
$ws.Range("N136").Value = -12028.2
$ws.Range("H138").Value = 53905.25
$ws.Range("J138").Value = 49970.668
$ws.Range("L138").Value = 49970.668
$ws.Range("N138").Value = -60250.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 74.75
$ws.Range("I40").Value = 56.857143
$ws.Range("J40").Value = 200
$ws.Range("K40").Value = 227.428572
$ws.Range("L40").Value = 800
$ws.Range("M40").Value = -158.428572
$ws.Range("N40").Value = -938
$ws.Range("H137").Value = 6901.4
$ws.Range("J137").Value = 10033.143
$ws.Range("L137").Value = 30099.429
$ws.Range("N137").Value = -40299.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 181344.64
$ws.Range("I70").Value = 194986.53
$ws.Range("K70").Value = 194986.53
$ws.Range("M70").Value = -194716.53
$ws.Range("H73").Value = 181344.64
$ws.Range("I73").Value = 194986.53
$ws.Range("K73").Value = 194986.53
$ws.Range("M73").Value = -194050.53
$ws.Range("H97").Value = 1860.3077
$ws.Range("I97").Value = 1439.2222
$ws.Range("K97").Value = 1439.2222
$ws.Range("M97").Value = -943.2221999999999
$ws.Range("H110").Value = 80918.914
$ws.Range("J110").Value = 80918.914
$ws.Range("L110").Value = 80918.914
$ws.Range("N110").Value = -89098.914
$ws.Range("H135").Value = 94991.664
$ws.Range("J135").Value = 94991.664
$ws.Range("L135").Value = 94991.664
$ws.Range("N135").Value = -105131.664
$ws.Range("H140").Value = 79329.89
$ws.Range("J140").Value = 90567
$ws.Range("L140").Value = 90567
$ws.Range("N140").Value = -100927

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4642
$ws.Range("I82").Value = 5845.5
$ws.Range("J82").Value = 3197.8
$ws.Range("K82").Value = 5845.5
$ws.Range("L82").Value = 3197.8
$ws.Range("M82").Value = -5484.5
$ws.Range("N82").Value = -3919.8
$ws.Range("H85").Value = 4642
$ws.Range("I85").Value = 5845.5
$ws.Range("J85").Value = 3197.8
$ws.Range("K85").Value = 5845.5
$ws.Range("L85").Value = 3197.8
$ws.Range("M85").Value = -4597.5
$ws.Range("N85").Value = -5693.8
$ws.Range("H100").Value = 11372.292
$ws.Range("I100").Value = 12441.75
$ws.Range("J100").Value = 9233.375
$ws.Range("K100").Value = 12441.75
$ws.Range("L100").Value = 9233.375
$ws.Range("M100").Value = -11900.75
$ws.Range("N100").Value = -10315.375
$ws.Range("H118").Value = 57779.2
$ws.Range("J118").Value = 57779.2
$ws.Range("L118").Value = 57779.2
$ws.Range("N118").Value = -61093.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 12432.8
$ws.Range("J41").Value = 12755
$ws.Range("L41").Value = 12755
$ws.Range("N41").Value = -13535
$ws.Range("H132").Value = 1404402.9
$ws.Range("I132").Value = 1697.2727
$ws.Range("J132").Value = 4833239
$ws.Range("K132").Value = 5091.8181
$ws.Range("L132").Value = 14499717
$ws.Range("M132").Value = -2561.8181
$ws.Range("N132").Value = -14504777
$ws.Range("H136").Value = 1343.6364
$ws.Range("I136").Value = 996.6667
$ws.Range("J136").Value = 2905
$ws.Range("K136").Value = 2990.0001
$ws.Range("L136").Value = 8715
$ws.Range("M136").Value = -440.0001000000002
$ws.Range("N136").Value = -13815

